# Tripadvisor New Orleans shard 206 update:
#  1. Reorder worksheets so "review_info" comes before "hotel_info".
#  2. Insert a new "State" column into hotel_info (between Hotel_Name and City)
#     and populate it with "Louisiana" for the existing hotel row.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the "State" column into hotel_info -------------------------
$hotelSheet = $wb.Worksheets.Item("hotel_info")
$hotelSheet.Range("C1").EntireColumn.Insert()
$hotelSheet.Range("C1").Value = "State"
$hotelSheet.Range("C2").Value = "Louisiana"

# --- 2. Move review_info in front of hotel_info ----------------------------
$reviewSheet = $wb.Worksheets.Item("review_info")
$reviewSheet.Move($wb.Worksheets.Item(1))
